# Atualizacao da resolucao do Desafio 02
# Adds a "Nome" column (C) with the product name, fills in the
# already-existing rows with an empty value in that column, and appends a
# new row (5) with the full set of data including the product name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell C1 ("Nome") -----------------------------------------
# Copy A1 first so C1 picks up the exact same (bold/centered/bordered)
# header style already used by A1 and B1, then overwrite the text.
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Nome"

# --- Existing rows 2-4: add an empty "Nome" cell -------------------------
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = ""
$ws.Range("C2").ClearFormats()

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = ""
$ws.Range("C3").ClearFormats()

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = ""
$ws.Range("C4").ClearFormats()

# --- New row 5 with full data, including the product name ---------------
# Force column A to stay plain text (otherwise "2024-08-22" would be
# auto-converted into a date serial number), then drop the format so the
# cell ends up with no explicit style, matching the rest of column A.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-08-22"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 4310.04

$ws.Range("C5").Value = "Apple iPhone 14 (128 GB) – Meia-Noite"
